# Update algorithm result values (commit: "Update Name of Algo")
# Applies the numeric value updates from the diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.344399999999999
$ws.Range("D3").Value = -7.014299999999999
$ws.Range("E19").Value = 16.2777
$ws.Range("B21").Value = 9.293800000000005
$ws.Range("B23").Value = 9.024200000000008
$ws.Range("D24").Value = -7.379500000000008
$ws.Range("E24").Value = 16.63540000000001
$ws.Range("B25").Value = 6.456700000000001
$ws.Range("C27").Value = -12.4369
$ws.Range("E30").Value = 15.5276
$ws.Range("C31").Value = -13.26799999999999
$ws.Range("E31").Value = 16.23310000000001
$ws.Range("E33").Value = 17.02190000000002
$ws.Range("C39").Value = -12.62070000000001
$ws.Range("C48").Value = -11.8043
$ws.Range("C51").Value = -11.7462
$ws.Range("C52").Value = -11.0246
$ws.Range("B53").Value = 5.420799999999998
$ws.Range("C55").Value = -13.7926
$ws.Range("E55").Value = 16.51070000000001
$ws.Range("C56").Value = -12.77979999999999
$ws.Range("B57").Value = 4.781899999999996
$ws.Range("C57").Value = -13.7887
$ws.Range("D57").Value = -8.9229
$ws.Range("B59").Value = 4.692799999999997
$ws.Range("D61").Value = -7.748999999999998
$ws.Range("E65").Value = 17.07750000000001
$ws.Range("B69").Value = 5.399899999999994
$ws.Range("D70").Value = -8.048799999999998
$ws.Range("E70").Value = 16.5526
$ws.Range("C73").Value = -12.45300000000001
$ws.Range("E75").Value = 16.51100000000001
$ws.Range("B79").Value = 9.5036
$ws.Range("B83").Value = 5.038599999999997
$ws.Range("E83").Value = 16.82
$ws.Range("D86").Value = -7.599499999999995
$ws.Range("C89").Value = -11.03240000000001
$ws.Range("C90").Value = -12.4138
$ws.Range("B93").Value = 5.5039
$ws.Range("E96").Value = 15.7693
$ws.Range("E97").Value = 16.8749
$ws.Range("D98").Value = -9.03779999999999
$ws.Range("D100").Value = -8.470599999999999
$ws.Range("D102").Value = -7.484299999999996
